$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (imie/nazwisko/login/mail/funkcja style record), matching
# the "getery imie nazwisko mail funkcja" commit: an extra user entry whose
# fields are unknown / 123 / 123@mail / brak.
$ws.Range("A3").Value = "unknown"
$ws.Range("B3").Value = "unknown"

# C3 / F3 hold the purely-numeric-looking text "123". A plain .Value
# assignment would be auto-coerced to a number by the engine, so force a
# text number-format on a scratch cell, paste the literal text in as a
# shared string there, then copy/paste-values it into place and restore the
# scratch cell - this keeps C3/F3 free of any leftover per-cell style index
# (matching how the rest of the sheet has no "s" attribute on its cells).
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "123"
$scratch.Copy()
$ws.Range("C3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("F3").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Range("D3").Value = "123@mail"
$ws.Range("E3").Value = "brak"
